$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 6.495
$ws.Range("A3").Value = -21.611
$ws.Range("B5").Value = 6.284000000000001
$ws.Range("C5").Value = -12.404
$ws.Range("E7").Value = 13.078
$ws.Range("C9").Value = -12.09
$ws.Range("C11").Value = -12.578
$ws.Range("E11").Value = 12.863
$ws.Range("A14").Value = -20.945
$ws.Range("A16").Value = -21.304
$ws.Range("B16").Value = 6.636999999999999
$ws.Range("C17").Value = -11.977
$ws.Range("E19").Value = 12.931
$ws.Range("A21").Value = -21.374
$ws.Range("C21").Value = -12.057
$ws.Range("E21").Value = 13.123
$ws.Range("A23").Value = -21.709
$ws.Range("A25").Value = -22.078
